$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fix spelling / wording of the header labels for the summary table (I8:N8)
$ws.Range("J8").Value = "total"
$ws.Range("K8").Value = "equal"
$ws.Range("L8").Value = "varied"
$ws.Range("M8").Value = "equal %"
$ws.Range("N8").Value = "varied %"

# The "Total" row label (I12) uses the same wording fix
$ws.Range("I12").Value = "total"

# Clear the (invisible / no-op) extra number-format & border styling that was
# previously applied to J9:L11 and M9:N11, reverting them to the plain
# default / 0.0-number-format styles (Excel drops the now-unused xf records
# on save).
$ws.Range("J9:L11").ClearFormats()
$ws.Range("M9:M11").NumberFormat = "0.0"
$ws.Range("N9:N11").NumberFormat = "0.0"

# Update the last active cell selection
$ws.Range("K18").Select()
